$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 12 (shifts rows 12-27 down to 13-28,
# and merged-cell ranges below shift automatically).
$ws.Rows("12:12").Insert()

# Row 12 becomes the new "Mean Live Births (SD)" summary row, matching
# the layout of row 2 ("Mean Age (SD)") and row 21 (old row 20).
$ws.Rows("12:12").RowHeight = 19

$ws.Range("A12").Value = "Mean Live Births (SD)"
$ws.Range("B12").Value = "2.36 (2.07)"
$ws.Range("C12").Value = "2.34 (2.11)"
$ws.Range("D12").Value = "t17480= 0.624"
$ws.Range("E12").Value = 0.532

# Subscript the "17480" portion of the test-statistic label in D12.
$dChars = $ws.Range("D12").Characters(2, 5)
$dChars.Font.Subscript = $true
$dChars.Font.Name = "Calibri (Body)"

# Formatting to match the sibling header rows (row 2 / row 21).
$ws.Range("A12").HorizontalAlignment = -4108
$ws.Range("A12").VerticalAlignment = -4108

$ws.Range("B12:D12").HorizontalAlignment = -4108

$ws.Range("E12").HorizontalAlignment = -4108
$ws.Range("E12").VerticalAlignment = -4108

# Restore the selection to A12 (matches the saved workbook view).
$ws.Range("A12").Select()
